$wb = $excel.ActiveWorkbook

# New values for column F, rows 4-14 (same for both sheets except row 10)
$newValues = @{
    4  = 74
    5  = 537
    6  = 7477
    7  = 469
    8  = 188
    9  = 1068
    11 = 17
    12 = 161
    13 = 195
    14 = 699
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}

# Row 10 differs between the two sheets (533->563 vs 534->564)
$wb.Worksheets.Item("展览").Range("F10").Value = 563
$wb.Worksheets.Item("全部类型").Range("F10").Value = 564
